$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 97

# Column A needs the text "01-08-2021" stored as a plain string (not auto-converted
# to a date serial). Typing it directly via .Value triggers Excel's date
# autodetection, so build it as a text formula first, then convert the formula
# to a static value in place (Copy + PasteSpecial values) which keeps the
# default "General" style untouched.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Formula = "=""01-08-2021"""
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = 202
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 0
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 202
$ws.Cells.Item($newRow, 11).Value = 0
